$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 10.23696233333333
$ws.Range("H2").Value = 30.710887
$ws.Range("I2").Value = 0.7155854078011842
$ws.Range("J2").Value = 0.7155854078011841
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2262196666666666
$ws.Range("N2").Value = 0.6786589999999999
$ws.Range("O2").Value = 0.03145179203784564
$ws.Range("P2").Value = 0.03145179203784564
$ws.Range("Q2").Value = 2.315802206725888
$ws.Range("R2").Value = 20.842219860533
$ws.Range("S2").Value = 0.02250644343147981
$ws.Range("T2").Value = 0.02250644343147981

# Row 3
$ws.Range("G3").Value = 10.23696233333333
$ws.Range("H3").Value = 30.710887
$ws.Range("I3").Value = 0.7155854078011842
$ws.Range("J3").Value = 0.7155854078011841
$ws.Range("O3").Value = 0.9636438974901603
$ws.Range("P3").Value = 0.9636438974901604
$ws.Range("Q3").Value = 70.95330725894333
$ws.Range("R3").Value = 638.57976533049
$ws.Range("S3").Value = 0.6895695113606189
$ws.Range("T3").Value = 0.6895695113606188

# Row 4
$ws.Range("G4").Value = 10.23696233333333
$ws.Range("H4").Value = 30.710887
$ws.Range("I4").Value = 0.7155854078011842
$ws.Range("J4").Value = 0.7155854078011841
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03527466666666667
$ws.Range("N4").Value = 0.105824
$ws.Range("O4").Value = 0.004904310471994002
$ws.Range("P4").Value = 0.004904310471994003
$ws.Range("Q4").Value = 0.3611054339875556
$ws.Range("R4").Value = 3.249948905888
$ws.Range("S4").Value = 0.003509453009085446
$ws.Range("T4").Value = 0.003509453009085446

# Row 5
$ws.Range("H5").Value = 4.509246
$ws.Range("I5").Value = 0.1050686239634127
$ws.Range("J5").Value = 0.1050686239634127
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2262196666666666
$ws.Range("N5").Value = 0.6786589999999999
$ws.Range("O5").Value = 0.03145179203784564
$ws.Range("P5").Value = 0.03145179203784564
$ws.Range("Q5").Value = 0.3400267090126666
$ws.Range("R5").Value = 3.060240381114
$ws.Range("S5").Value = 0.00330459651059986
$ws.Range("T5").Value = 0.00330459651059986

# Row 6
$ws.Range("H6").Value = 4.509246
$ws.Range("I6").Value = 0.1050686239634127
$ws.Range("J6").Value = 0.1050686239634127
$ws.Range("O6").Value = 0.9636438974901603
$ws.Range("P6").Value = 0.9636438974901604
$ws.Range("R6").Value = 93.76196957442001
$ws.Range("S6").Value = 0.101248738300031
$ws.Range("T6").Value = 0.101248738300031

# Row 7
$ws.Range("H7").Value = 4.509246
$ws.Range("I7").Value = 0.1050686239634127
$ws.Range("J7").Value = 0.1050686239634127
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.03527466666666667
$ws.Range("N7").Value = 0.105824
$ws.Range("O7").Value = 0.004904310471994002
$ws.Range("P7").Value = 0.004904310471994003
$ws.Range("Q7").Value = 0.05302071652266667
$ws.Range("R7").Value = 0.477186448704
$ws.Range("S7").Value = 0.0005152891527817647
$ws.Range("T7").Value = 0.0005152891527817647

# Row 8
$ws.Range("G8").Value = 1.544659666666667
$ws.Range("H8").Value = 4.633979
$ws.Range("I8").Value = 0.107974991163789
$ws.Range("J8").Value = 0.107974991163789
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.2262196666666666
$ws.Range("N8").Value = 0.6786589999999999
$ws.Range("O8").Value = 0.03145179203784564
$ws.Range("P8").Value = 0.03145179203784564
$ws.Range("Q8").Value = 0.3494323949067777
$ws.Range("R8").Value = 3.144891554161
$ws.Range("S8").Value = 0.003396006967371714
$ws.Range("T8").Value = 0.003396006967371714

# Row 9
$ws.Range("G9").Value = 1.544659666666667
$ws.Range("H9").Value = 4.633979
$ws.Range("I9").Value = 0.107974991163789
$ws.Range("J9").Value = 0.107974991163789
$ws.Range("O9").Value = 0.9636438974901603
$ws.Range("P9").Value = 0.9636438974901604
$ws.Range("Q9").Value = 10.70617516903667
$ws.Range("R9").Value = 96.35557652132999
$ws.Range("S9").Value = 0.1040494413165393
$ws.Range("T9").Value = 0.1040494413165393

# Row 10
$ws.Range("G10").Value = 1.544659666666667
$ws.Range("H10").Value = 4.633979
$ws.Range("I10").Value = 0.107974991163789
$ws.Range("J10").Value = 0.107974991163789
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.03527466666666667
$ws.Range("N10").Value = 0.105824
$ws.Range("O10").Value = 0.004904310471994002
$ws.Range("P10").Value = 0.004904310471994003
$ws.Range("Q10").Value = 0.05448735485511112
$ws.Range("R10").Value = 0.490386193696
$ws.Range("S10").Value = 0.0005295428798780305
$ws.Range("T10").Value = 0.0005295428798780305

# Row 11
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.04013633333333334
$ws.Range("H11").Value = 0.120409
$ws.Range("I11").Value = 0.002805614939351403
$ws.Range("J11").Value = 0.002805614939351403
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.2262196666666666
$ws.Range("N11").Value = 0.6786589999999999
$ws.Range("O11").Value = 0.03145179203784564
$ws.Range("P11").Value = 0.03145179203784564
$ws.Range("Q11").Value = 0.009079627947888889
$ws.Range("R11").Value = 0.081716651531
$ws.Range("S11").Value = 0.00008824161761075325
$ws.Range("T11").Value = 0.00008824161761075324

# Row 12
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.04013633333333334
$ws.Range("H12").Value = 0.120409
$ws.Range("I12").Value = 0.002805614939351403
$ws.Range("J12").Value = 0.002805614939351403
$ws.Range("O12").Value = 0.9636438974901603
$ws.Range("P12").Value = 0.9636438974901604
$ws.Range("Q12").Value = 0.2781885386033334
$ws.Range("R12").Value = 2.50369684743
$ws.Range("S12").Value = 0.002703613715013206
$ws.Range("T12").Value = 0.002703613715013206

# Row 13
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.04013633333333334
$ws.Range("H13").Value = 0.120409
$ws.Range("I13").Value = 0.002805614939351403
$ws.Range("J13").Value = 0.002805614939351403
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.03527466666666667
$ws.Range("N13").Value = 0.105824
$ws.Range("O13").Value = 0.004904310471994002
$ws.Range("P13").Value = 0.004904310471994003
$ws.Range("Q13").Value = 0.001415795779555556
$ws.Range("R13").Value = 0.012742162016
$ws.Range("S13").Value = 0.0000137596067274439
$ws.Range("T13").Value = 0.00001375960672744391

# Row 14
$ws.Range("G14").Value = 0.5501683333333333
$ws.Range("H14").Value = 1.650505
$ws.Range("I14").Value = 0.03845793491744127
$ws.Range("J14").Value = 0.03845793491744127
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.2262196666666666
$ws.Range("N14").Value = 0.6786589999999999
$ws.Range("O14").Value = 0.03145179203784564
$ws.Range("P14").Value = 0.03145179203784564
$ws.Range("Q14").Value = 0.1244588969772222
$ws.Range("R14").Value = 1.120130072795
$ws.Range("S14").Value = 0.001209570971228365
$ws.Range("T14").Value = 0.001209570971228365

# Row 15
$ws.Range("G15").Value = 0.5501683333333333
$ws.Range("H15").Value = 1.650505
$ws.Range("I15").Value = 0.03845793491744127
$ws.Range("J15").Value = 0.03845793491744127
$ws.Range("O15").Value = 0.9636438974901603
$ws.Range("P15").Value = 0.9636438974901604
$ws.Range("Q15").Value = 3.813266233483333
$ws.Range("R15").Value = 34.31939610134999
$ws.Range("S15").Value = 0.03705975429326604
$ws.Range("T15").Value = 0.03705975429326604

# Row 16
$ws.Range("G16").Value = 0.5501683333333333
$ws.Range("H16").Value = 1.650505
$ws.Range("I16").Value = 0.03845793491744127
$ws.Range("J16").Value = 0.03845793491744127
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.03527466666666667
$ws.Range("N16").Value = 0.105824
$ws.Range("O16").Value = 0.004904310471994002
$ws.Range("P16").Value = 0.004904310471994003
$ws.Range("Q16").Value = 0.01940700456888889
$ws.Range("R16").Value = 0.17466304112
$ws.Range("S16").Value = 0.000188609652946871
$ws.Range("T16").Value = 0.0001886096529468711

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.4307083333333333
$ws.Range("H17").Value = 1.292125
$ws.Range("I17").Value = 0.03010742721482141
$ws.Range("J17").Value = 0.03010742721482141
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.2262196666666666
$ws.Range("N17").Value = 0.6786589999999999
$ws.Range("O17").Value = 0.03145179203784564
$ws.Range("P17").Value = 0.03145179203784564
$ws.Range("Q17").Value = 0.0974346955972222
$ws.Range("R17").Value = 0.8769122603749998
$ws.Range("S17").Value = 0.0009469325395551373
$ws.Range("T17").Value = 0.0009469325395551373

# Row 18
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 0.6666666666666666
$ws.Range("G18").Value = 0.4307083333333333
$ws.Range("H18").Value = 1.292125
$ws.Range("I18").Value = 0.03010742721482141
$ws.Range("J18").Value = 0.03010742721482141
$ws.Range("O18").Value = 0.9636438974901603
$ws.Range("P18").Value = 0.9636438974901604
$ws.Range("Q18").Value = 2.985278222083333
$ws.Range("R18").Value = 26.86750399875
$ws.Range("S18").Value = 0.02901283850469183
$ws.Range("T18").Value = 0.02901283850469183

# Row 19
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 0.6666666666666666
$ws.Range("G19").Value = 0.4307083333333333
$ws.Range("H19").Value = 1.292125
$ws.Range("I19").Value = 0.03010742721482141
$ws.Range("J19").Value = 0.03010742721482141
$ws.Range("K19").Value = 1
$ws.Range("L19").Value = 0.3333333333333333
$ws.Range("M19").Value = 0.03527466666666667
$ws.Range("N19").Value = 0.105824
$ws.Range("O19").Value = 0.004904310471994002
$ws.Range("P19").Value = 0.004904310471994003
$ws.Range("Q19").Value = 0.01519309288888889
$ws.Range("R19").Value = 0.136737836
$ws.Range("S19").Value = 0.0001476561705744459
$ws.Range("T19").Value = 0.0001476561705744459
